$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns D, M, N, O, P, S per row (row 4 unchanged)
$values = @{
    2  = @(44475, 240, 11000, 12000, 11500, 5750)
    3  = @(44454, 160, 12000, 13000, 12500, 6250)
    5  = @(44482, 240, 10000, 11000, 10500, 5250)
    6  = @(44517, 400, 5500,  6000,  5750,  2875)
    7  = @(44461, 200, 11000, 12000, 11500, 5750)
    8  = @(44489, 160, 9500,  10000, 9750,  4875)
    9  = @(44497, 500, 9000,  10000, 9500,  4750)
    10 = @(44490, 400, 9500,  10000, 9750,  4875)
}

foreach ($row in $values.Keys) {
    $v = $values[$row]
    $ws.Range("D$row").Value = $v[0]
    $ws.Range("M$row").Value = $v[1]
    $ws.Range("N$row").Value = $v[2]
    $ws.Range("O$row").Value = $v[3]
    $ws.Range("P$row").Value = $v[4]
    $ws.Range("S$row").Value = $v[5]
}
